# Convert the clueBoard sheet into a clean .csv-ready grid:
# drop the helper row-index column X (X1:X22) and the helper header
# row 23 (A23:W23) that were used for import bookkeeping, then leave
# the selection on the old header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the auxiliary row-index values in column X (rows 1-22).
$ws.Range("X1:X22").ClearContents()

# Remove the auxiliary column-index header row (row 23).
$ws.Range("A23:W23").ClearContents()

# Match the saved selection state.
$ws.Range("A23:W23").Select()
